$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 66 ("WHD110") and push the existing
# WHD110/WHD120/BMXWT rows down by one, then fill in the new variable name.
$ws.Rows.Item(66).Insert()
$ws.Range("A66").Value = "HSD010"

# Match the saved selection/active cell from the edit.
$ws.Range("I58").Select() | Out-Null
